# Apply automated update values (2025-10-13 10:30:09)
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M12").Value = 4543.49

# --- Sheet: VENTA MENSUAL ---
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F12").Value = 4543.49
$wsVentaMensual.Range("F26").Value = 9672.959999999999

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 11064.99
$wsCumplimiento.Range("E12").Value = 16889.99
$wsCumplimiento.Range("F12").Value = 0.3958146276620481

$wsCumplimiento.Range("D14").Value = 9672.959999999999
$wsCumplimiento.Range("E14").Value = 32530.42110009469
$wsCumplimiento.Range("F14").Value = 0.2291986980156502
